$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 8547612
$ws.Range("I53").Value = 20833780
$ws.Range("J53").Value = 712.1739
$ws.Range("K53").Value = 20833780
$ws.Range("L53").Value = 712.1739
$ws.Range("M53").Value = -20833143
$ws.Range("N53").Value = -1986.1739

$ws.Range("H98").Value = 1559.4038
$ws.Range("I98").Value = 1195.4147
$ws.Range("K98").Value = 1195.4147
$ws.Range("M98").Value = 302.5853

$ws.Range("H107").Value = 32391.938
$ws.Range("I107").Value = 42498.293
$ws.Range("J107").Value = 2072.875
$ws.Range("K107").Value = 42498.293
$ws.Range("L107").Value = 2072.875
$ws.Range("M107").Value = -40578.293
$ws.Range("N107").Value = -5912.875

$ws.Range("H122").Value = 1559.4038
$ws.Range("I122").Value = 1195.4147
$ws.Range("K122").Value = 3586.2441
$ws.Range("M122").Value = -1136.2441

$ws.Range("H131").Value = 1970.0938
$ws.Range("I131").Value = 1666.5186
$ws.Range("J131").Value = 3609.4
$ws.Range("K131").Value = 4999.5558
$ws.Range("L131").Value = 10828.2
$ws.Range("M131").Value = 40.44419999999991
$ws.Range("N131").Value = -20908.2

$ws.Range("H132").Value = 1472.1729
$ws.Range("I132").Value = 1558.625
$ws.Range("K132").Value = 4675.875
$ws.Range("M132").Value = -2145.875

$ws.Range("H137").Value = 5271.933
$ws.Range("I137").Value = 3786
$ws.Range("K137").Value = 11358
$ws.Range("M137").Value = -8808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 117.55556
$ws.Range("I5").Value = 119.75
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 119.75
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -7.75
$ws.Range("N5").Value = -324

$ws.Range("H32").Value = 2231.76
$ws.Range("I32").Value = 1836.8969
$ws.Range("K32").Value = 1836.8969
$ws.Range("M32").Value = -1549.8969

$ws.Range("H45").Value = 3874.7693
$ws.Range("I45").Value = 2767.2778
$ws.Range("J45").Value = 6366.625
$ws.Range("K45").Value = 2767.2778
$ws.Range("L45").Value = 6366.625
$ws.Range("M45").Value = -2390.2778
$ws.Range("N45").Value = -7120.625

$ws.Range("H61").Value = 2870.1853
$ws.Range("I61").Value = 1737.7858
$ws.Range("K61").Value = 1737.7858
$ws.Range("M61").Value = -1525.7858

$ws.Range("H132").Value = 6010.0312
$ws.Range("J132").Value = 7913.3335
$ws.Range("L132").Value = 23740.0005
$ws.Range("N132").Value = -28800.0005

$ws.Range("H136").Value = 2870.1853
$ws.Range("I136").Value = 1737.7858
$ws.Range("K136").Value = 5213.357400000001
$ws.Range("M136").Value = -2663.357400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 117.55556
$ws.Range("I4").Value = 119.75
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 119.75
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -4.75
$ws.Range("N4").Value = -330

$ws.Range("H22").Value = 262.5
$ws.Range("I22").Value = 325
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 325
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = -152
$ws.Range("N22").Value = -546

$ws.Range("H99").Value = 5099.0454
$ws.Range("I99").Value = 5182.4443
$ws.Range("J99").Value = 4723.75
$ws.Range("K99").Value = 5182.4443
$ws.Range("L99").Value = 4723.75
$ws.Range("M99").Value = -3684.4443
$ws.Range("N99").Value = -7719.75

$ws.Range("H133").Value = 74987.5
$ws.Range("J133").Value = 74987.5
$ws.Range("L133").Value = 74987.5
$ws.Range("N133").Value = -85107.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 772503
$ws.Range("I31").Value = 835211.5600000001
$ws.Range("K31").Value = 835211.5600000001
$ws.Range("M31").Value = -834916.5600000001

$ws.Range("H34").Value = 772503
$ws.Range("I34").Value = 835211.5600000001
$ws.Range("K34").Value = 835211.5600000001
$ws.Range("M34").Value = -835009.5600000001

$ws.Range("H132").Value = 2982.484
$ws.Range("I132").Value = 2132.75
$ws.Range("J132").Value = 5895.857
$ws.Range("K132").Value = 6398.25
$ws.Range("L132").Value = 17687.571
$ws.Range("M132").Value = -3868.25
$ws.Range("N132").Value = -22747.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3441.5881
$ws.Range("I2").Value = 7006.625
$ws.Range("K2").Value = 42039.75
$ws.Range("M2").Value = -41926.75

$ws.Range("H14").Value = 388.14285
$ws.Range("I14").Value = 388.14285
$ws.Range("K14").Value = 1164.42855
$ws.Range("M14").Value = -991.4285500000001

$ws.Range("H23").Value = 187.23077
$ws.Range("I23").Value = 280.85715
$ws.Range("J23").Value = 78
$ws.Range("K23").Value = 842.5714499999999
$ws.Range("L23").Value = 234
$ws.Range("M23").Value = -607.5714499999999
$ws.Range("N23").Value = -704

$ws.Range("H86").Value = 291.5
$ws.Range("J86").Value = 291.5
$ws.Range("L86").Value = 874.5
$ws.Range("N86").Value = -3246.5

$ws.Range("H89").Value = 291.5
$ws.Range("J89").Value = 291.5
$ws.Range("L89").Value = 2623.5
$ws.Range("N89").Value = -14479.5

$ws.Range("H104").Value = 15000
$ws.Range("J104").Value = 15000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -50242

$ws.Range("H122").Value = 48693.047
$ws.Range("J122").Value = 84730.664
$ws.Range("L122").Value = 762575.976
$ws.Range("N122").Value = -767475.976

$ws.Range("H132").Value = 4169.0835
$ws.Range("I132").Value = 895.6667
$ws.Range("J132").Value = 5260.222
$ws.Range("K132").Value = 8061.0003
$ws.Range("L132").Value = 47341.998
$ws.Range("M132").Value = -5531.0003
$ws.Range("N132").Value = -52401.998

$ws.Range("H134").Value = 1807.125
$ws.Range("I134").Value = 1636.7142
$ws.Range("K134").Value = 4910.142599999999
$ws.Range("M134").Value = 159.8574000000008

$ws.Range("H136").Value = 8749.75
$ws.Range("I136").Value = 7666.3335
$ws.Range("K136").Value = 22999.0005
$ws.Range("M136").Value = -17899.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 500072500
$ws.Range("J18").Value = 500072500
$ws.Range("L18").Value = 500072500
$ws.Range("N18").Value = -500073086

$ws.Range("H99").Value = 5299.8
$ws.Range("I99").Value = 4999
$ws.Range("J99").Value = 5375
$ws.Range("K99").Value = 4999
$ws.Range("L99").Value = 5375
$ws.Range("M99").Value = -2753
$ws.Range("N99").Value = -9867

$ws.Range("H132").Value = 189242.36
$ws.Range("I132").Value = 224646.69
$ws.Range("J132").Value = 75442.71000000001
$ws.Range("K132").Value = 673940.0700000001
$ws.Range("L132").Value = 226328.13
$ws.Range("M132").Value = -671410.0700000001
$ws.Range("N132").Value = -231388.13

$ws.Range("H141").Value = 28900
$ws.Range("J141").Value = 28900
$ws.Range("L141").Value = 28900
$ws.Range("N141").Value = -39260

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6711.222
$ws.Range("I22").Value = 5433.5
$ws.Range("K22").Value = 5433.5
$ws.Range("M22").Value = -5138.5

$ws.Range("H27").Value = 6711.222
$ws.Range("I27").Value = 5433.5
$ws.Range("K27").Value = 5433.5
$ws.Range("M27").Value = -5326.5

$ws.Range("H46").Value = 5125.25
$ws.Range("J46").Value = 6002
$ws.Range("L46").Value = 6002
$ws.Range("N46").Value = -6378

$ws.Range("H105").Value = 11500
$ws.Range("J105").Value = 11500
$ws.Range("L105").Value = 11500
$ws.Range("N105").Value = -18488

$ws.Range("H132").Value = 3665.054
$ws.Range("I132").Value = 2948.6
$ws.Range("J132").Value = 6735.5713
$ws.Range("K132").Value = 8845.799999999999
$ws.Range("L132").Value = 20206.7139
$ws.Range("M132").Value = -6315.799999999999
$ws.Range("N132").Value = -25266.7139

$ws.Range("H136").Value = 3384.9656
$ws.Range("I136").Value = 2731.2942
$ws.Range("K136").Value = 8193.882599999999
$ws.Range("M136").Value = -5643.882599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7398
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 7398
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 7398
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -9270

$ws.Range("H77").Value = 7398
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 7398
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 22194
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -31554

$ws.Range("H81").Value = 19523.375
$ws.Range("I81").Value = 2065.3333
$ws.Range("K81").Value = 4130.6666
$ws.Range("M81").Value = -3069.6666

$ws.Range("H84").Value = 19523.375
$ws.Range("I84").Value = 2065.3333
$ws.Range("K84").Value = 20653.333
$ws.Range("M84").Value = -15349.333

$ws.Range("H100").Value = 1167
$ws.Range("I100").Value = 1001
$ws.Range("K100").Value = 2002
$ws.Range("M100").Value = -1461

$ws.Range("H132").Value = 19096.268
$ws.Range("I132").Value = 1661.5227
$ws.Range("K132").Value = 4984.5681
$ws.Range("M132").Value = -2454.5681

$ws.Range("H136").Value = 55104.957
$ws.Range("I136").Value = 11041.023
$ws.Range("K136").Value = 33123.069
$ws.Range("M136").Value = -30573.069

$ws.Range("H140").Value = 78270.57000000001
$ws.Range("J140").Value = 78270.57000000001
$ws.Range("L140").Value = 78270.57000000001
$ws.Range("N140").Value = -88630.57000000001
